# Update TPM-derived NATMI ligand-receptor metrics (Col2a1-Itga2, YoungD7)
# with recomputed values reflecting the new TPM-based ligand/receptor
# expressing-cell counts and downstream specificity/weight scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.107333666666666
$ws.Range("N2").Value = 21.322001
$ws.Range("O2").Value = 0.7373665550576455
$ws.Range("P2").Value = 0.7373665550576454
$ws.Range("Q2").Value = 1.124363602288111
$ws.Range("R2").Value = 10.119272420593
$ws.Range("S2").Value = 0.1032804299525234
$ws.Range("T2").Value = 0.1032804299525234

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("O3").Value = 0.1688878844614928
$ws.Range("P3").Value = 0.1688878844614928
$ws.Range("Q3").Value = 0.2575264484854444
$ws.Range("R3").Value = 2.317738036369
$ws.Range("S3").Value = 0.02365555258956848
$ws.Range("T3").Value = 0.02365555258956847

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1581976666666667
$ws.Range("H4").Value = 0.474593
$ws.Range("I4").Value = 0.1400666049254827
$ws.Range("J4").Value = 0.1400666049254826
$ws.Range("M4").Value = 0.8135026666666666
$ws.Range("N4").Value = 2.440508
$ws.Range("O4").Value = 0.08439869112428164
$ws.Range("P4").Value = 0.08439869112428162
$ws.Range("Q4").Value = 0.1286942236937778
$ws.Range("R4").Value = 1.158248013244
$ws.Range("S4").Value = 0.0118214381259326
$ws.Range("T4").Value = 0.01182143812593259

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1581976666666667
$ws.Range("H5").Value = 0.474593
$ws.Range("I5").Value = 0.1400666049254827
$ws.Range("J5").Value = 0.1400666049254826
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09009266666666667
$ws.Range("N5").Value = 0.270278
$ws.Range("O5").Value = 0.009346869356580103
$ws.Range("P5").Value = 0.009346869356580103
$ws.Range("Q5").Value = 0.01425244965044445
$ws.Range("R5").Value = 0.128272046854
$ws.Range("S5").Value = 0.001309184257458206
$ws.Range("T5").Value = 0.001309184257458205

# Row 6
$ws.Range("I6").Value = 0.7029419733214338
$ws.Range("J6").Value = 0.7029419733214337
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.107333666666666
$ws.Range("N6").Value = 21.322001
$ws.Range("O6").Value = 0.7373665550576455
$ws.Range("P6").Value = 0.7373665550576454
$ws.Range("Q6").Value = 5.642760954644999
$ws.Range("R6").Value = 50.784848591805
$ws.Range("S6").Value = 0.518325901273449
$ws.Range("T6").Value = 0.5183259012734488

# Row 7
$ws.Range("I7").Value = 0.7029419733214338
$ws.Range("J7").Value = 0.7029419733214337
$ws.Range("O7").Value = 0.1688878844614928
$ws.Range("P7").Value = 0.1688878844614928
$ws.Range("S7").Value = 0.1187183827734441
$ws.Range("T7").Value = 0.118718382773444

# Row 8
$ws.Range("I8").Value = 0.7029419733214338
$ws.Range("J8").Value = 0.7029419733214337
$ws.Range("M8").Value = 0.8135026666666666
$ws.Range("N8").Value = 2.440508
$ws.Range("O8").Value = 0.08439869112428164
$ws.Range("P8").Value = 0.08439869112428162
$ws.Range("Q8").Value = 0.64586823966
$ws.Range("R8").Value = 5.81281415694
$ws.Range("S8").Value = 0.05932738248464871
$ws.Range("T8").Value = 0.0593273824846487

# Row 9
$ws.Range("I9").Value = 0.7029419733214338
$ws.Range("J9").Value = 0.7029419733214337
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09009266666666667
$ws.Range("N9").Value = 0.270278
$ws.Range("O9").Value = 0.009346869356580103
$ws.Range("P9").Value = 0.009346869356580103
$ws.Range("Q9").Value = 0.07152772131
$ws.Range("R9").Value = 0.6437494917900001
$ws.Range("S9").Value = 0.006570306789892058
$ws.Range("T9").Value = 0.006570306789892057

# Row 10
$ws.Range("G10").Value = 0.1199896666666667
$ws.Range("H10").Value = 0.359969
$ws.Range("I10").Value = 0.106237630366274
$ws.Range("J10").Value = 0.106237630366274
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.107333666666666
$ws.Range("N10").Value = 21.322001
$ws.Range("O10").Value = 0.7373665550576455
$ws.Range("P10").Value = 0.7373665550576454
$ws.Range("Q10").Value = 0.8528065975521111
$ws.Range("R10").Value = 7.675259377969
$ws.Range("S10").Value = 0.07833607552066696
$ws.Range("T10").Value = 0.07833607552066692

# Row 11
$ws.Range("G11").Value = 0.1199896666666667
$ws.Range("H11").Value = 0.359969
$ws.Range("I11").Value = 0.106237630366274
$ws.Range("J11").Value = 0.106237630366274
$ws.Range("O11").Value = 0.1688878844614928
$ws.Range("P11").Value = 0.1688878844614928
$ws.Range("Q11").Value = 0.1953284985974444
$ws.Range("R11").Value = 1.757956487377
$ws.Range("S11").Value = 0.01794224864276206
$ws.Range("T11").Value = 0.01794224864276205

# Row 12
$ws.Range("G12").Value = 0.1199896666666667
$ws.Range("H12").Value = 0.359969
$ws.Range("I12").Value = 0.106237630366274
$ws.Range("J12").Value = 0.106237630366274
$ws.Range("M12").Value = 0.8135026666666666
$ws.Range("N12").Value = 2.440508
$ws.Range("O12").Value = 0.08439869112428164
$ws.Range("P12").Value = 0.08439869112428162
$ws.Range("Q12").Value = 0.09761191380577776
$ws.Range("R12").Value = 0.8785072242519999
$ws.Range("S12").Value = 0.008966316951058761
$ws.Range("T12").Value = 0.008966316951058758

# Row 13
$ws.Range("G13").Value = 0.1199896666666667
$ws.Range("H13").Value = 0.359969
$ws.Range("I13").Value = 0.106237630366274
$ws.Range("J13").Value = 0.106237630366274
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09009266666666667
$ws.Range("N13").Value = 0.270278
$ws.Range("O13").Value = 0.009346869356580103
$ws.Range("P13").Value = 0.009346869356580103
$ws.Range("Q13").Value = 0.01081018904244444
$ws.Range("R13").Value = 0.097291701382
$ws.Range("S13").Value = 0.00099298925178621
$ws.Range("T13").Value = 0.0009929892517862098

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05732366666666666
$ws.Range("H14").Value = 0.171971
$ws.Range("I14").Value = 0.05075379138680971
$ws.Range("J14").Value = 0.05075379138680969
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.107333666666666
$ws.Range("N14").Value = 21.322001
$ws.Range("O14").Value = 0.7373665550576455
$ws.Range("P14").Value = 0.7373665550576454
$ws.Range("Q14").Value = 0.4074184259967777
$ws.Range("R14").Value = 3.666765833971
$ws.Range("S14").Value = 0.03742414831100627
$ws.Range("T14").Value = 0.03742414831100626

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05732366666666666
$ws.Range("H15").Value = 0.171971
$ws.Range("I15").Value = 0.05075379138680971
$ws.Range("J15").Value = 0.05075379138680969
$ws.Range("O15").Value = 0.1688878844614928
$ws.Range("P15").Value = 0.1688878844614928
$ws.Range("Q15").Value = 0.0933159167381111
$ws.Range("R15").Value = 0.8398432506429999
$ws.Range("S15").Value = 0.008571700455718229
$ws.Range("T15").Value = 0.008571700455718223

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05732366666666666
$ws.Range("H16").Value = 0.171971
$ws.Range("I16").Value = 0.05075379138680971
$ws.Range("J16").Value = 0.05075379138680969
$ws.Range("M16").Value = 0.8135026666666666
$ws.Range("N16").Value = 2.440508
$ws.Range("O16").Value = 0.08439869112428164
$ws.Range("P16").Value = 0.08439869112428162
$ws.Range("Q16").Value = 0.04663295569644443
$ws.Range("R16").Value = 0.419696601268
$ws.Range("S16").Value = 0.004283553562641578
$ws.Range("T16").Value = 0.004283553562641576

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05732366666666666
$ws.Range("H17").Value = 0.171971
$ws.Range("I17").Value = 0.05075379138680971
$ws.Range("J17").Value = 0.05075379138680969
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09009266666666667
$ws.Range("N17").Value = 0.270278
$ws.Range("O17").Value = 0.009346869356580103
$ws.Range("P17").Value = 0.009346869356580103
$ws.Range("Q17").Value = 0.00516444199311111
$ws.Range("R17").Value = 0.046479977938
$ws.Range("S17").Value = 0.0004743890574436308
$ws.Range("T17").Value = 0.0004743890574436307
